$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: simple in-place text replacement of a unique literal string that
# lives inside a single run (keeps that run's formatting untouched).
# ---------------------------------------------------------------------------
function Replace-Simple($old, $new) {
    $rng = $d.Content
    $ok = $rng.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $ok) {
        Write-Output "NOT FOUND (simple): $old"
    }
}

# ---------------------------------------------------------------------------
# Helper: replace every occurrence of a unique literal string (ReplaceAll).
# ---------------------------------------------------------------------------
function Replace-All($old, $new) {
    $rng = $d.Content
    $ok = $rng.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $ok) {
        Write-Output "NOT FOUND (all): $old"
    }
}

# ---------------------------------------------------------------------------
# Helper: the heading paragraphs store "<bold run1><non-bold run2>" where
# run2 used to start with some filler text and end in a colon. The new
# wording folds the filler into run1 (bold) and leaves run2 holding just the
# colon (still non-bold, its own run). We locate run1 by its old text,
# overwrite it in place (so that run keeps its own formatting), then
# overwrite the immediately following run2 text (also in place) so the two
# runs never collapse into one.
# ---------------------------------------------------------------------------
function Replace-HeadingColon($oldRun1, $newRun1, $oldRun2, $newRun2) {
    $rng = $d.Content
    $ok = $rng.Find.Execute($oldRun1, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        Write-Output "NOT FOUND (head run1): $oldRun1"
        return
    }
    $rng.Text = $newRun1
    $afterEnd = $rng.End
    $rng2 = $d.Range($afterEnd, $afterEnd + $oldRun2.Length)
    if ($rng2.Text -ne $oldRun2) {
        Write-Output "MISMATCH (head run2): expected [$oldRun2] got [$($rng2.Text)]"
    }
    $rng2.Text = $newRun2
}

# 1) Intro paragraph -------------------------------------------------------
Replace-Simple `
    "en función de las directrices del sector de TI e incluir los elementos adicionales que mencionó:" `
    "en función de las directrices del sector de TI y los elementos adicionales mencionados:"

# 2) "Evalúe la posición" / " de seguridad de red actual:" -----------------
Replace-HeadingColon `
    "Evalúe la posición" `
    "Evaluación de la posición de seguridad de red actual" `
    " de seguridad de red actual:" `
    ":"

Replace-Simple `
    "Realice una auditoría de seguridad para revisar las prácticas y directivas de seguridad de red actuales." `
    "Realización de una auditoría de seguridad para revisar las prácticas y directivas de seguridad de red actuales."

Replace-Simple `
    "Identifique las brechas o vulnerabilidades que deben abordarse." `
    "Identificación de las brechas o vulnerabilidades que deben abordarse."

Replace-Simple `
    "Revise los procedimientos recomendados del sector, como los recomendados por NIST, para garantizar el cumplimiento." `
    "Revisión de los procedimientos recomendados del sector, como los recomendados por NIST, para garantizar el cumplimiento."

# 3) "Seleccione el producto" / " de seguridad de red adecuado:" -----------
Replace-HeadingColon `
    "Seleccione el producto" `
    "Selección del producto de seguridad de red adecuado" `
    " de seguridad de red adecuado:" `
    ":"

Replace-Simple `
    "Investigue y evalúe diferentes productos de seguridad de red." `
    "Investigación y evaluación de diferentes productos de seguridad de red."

Replace-Simple `
    "Determine qué producto satisface mejor las necesidades de la organización, teniendo en cuenta factores como la compatibilidad con los sistemas existentes, la facilidad de uso y el costo." `
    "Determinación de qué producto satisface mejor las necesidades de la organización, teniendo en cuenta factores como la compatibilidad con los sistemas existentes, la facilidad de uso y el coste."

Replace-Simple `
    "En este caso, se ha seleccionado Contoso CipherGuard Sentinel X7 como producto de seguridad de red que se va a instalar." `
    "En este caso, se ha seleccionado Contoso CipherGuard Sentinel X7 como el producto de seguridad de red que se va a instalar."

# 4) "Desarrollar un plan" / " de implementación:" --------------------------
Replace-HeadingColon `
    "Desarrollar un plan" `
    "Desarrollo de un plan de implementación" `
    " de implementación:" `
    ":"

# 5) "Configure e instale el producto" / " de seguridad de red:" ------------
Replace-HeadingColon `
    "Configure e instale el producto" `
    "Configuración e instalación del producto de seguridad de red" `
    " de seguridad de red:" `
    ":"

# 6) "Pruebe y valide la implementación" (colon already its own run) -------
Replace-Simple `
    "Pruebe y valide la implementación" `
    "Prueba y validación de la implementación"

Replace-Simple `
    "Realice pruebas exhaustivas para asegurarse de que Contoso CipherGuard Sentinel X7 está configurado y funciona correctamente según lo previsto." `
    "Realice pruebas exhaustivas para asegurarse de que Contoso CipherGuard Sentinel X7 está configurado correctamente y funciona según lo previsto."

Replace-Simple `
    "Realice pruebas de penetración o examen de vulnerabilidades para identificar posibles debilidades." `
    "Realice pruebas de penetración o un examen de vulnerabilidades para identificar posibles debilidades."

# 7) "Entrenar usuarios y administradores" (colon already its own run) -----
Replace-Simple `
    "Entrenar usuarios y administradores" `
    "Formación de usuarios y administradores"

Replace-Simple `
    "Proporcione aprendizaje a los usuarios y administradores sobre cómo usar y mantener correctamente Contoso CipherGuard Sentinel X7." `
    "Proporcione formación a los usuarios y administradores sobre cómo usar y mantener correctamente Contoso CipherGuard Sentinel X7."

# occurs twice, same replacement both times
Replace-All `
    "Proporcione documentación, realice sesiones de entrenamiento o ofrezca soporte técnico continuo." `
    "Proporciona documentación, realiza sesiones de formación u ofrece soporte técnico continuo."

# 8) "Supervise y mantenga el producto" / " de seguridad de red:" ----------
Replace-HeadingColon `
    "Supervise y mantenga el producto" `
    "Supervisión y mantenimiento del producto de seguridad de red" `
    " de seguridad de red:" `
    ":"

Replace-Simple `
    "Supervise periódicamente Contoso CipherGuard Sentinel X7 para asegurarse de que funciona correctamente y proporciona el nivel de protección deseado." `
    "Supervise periódicamente Contoso CipherGuard Sentinel X7 para asegurarse de que funciona correctamente y proporcione el nivel de protección deseado."

Replace-Simple `
    "Desarrolle e implemente un programa de entrenamiento para asegurarse de que todos los usuarios y administradores estén entrenados correctamente sobre cómo usar y mantener Contoso CipherGuard Sentinel X7." `
    "Desarrolle e implemente un programa de formación para asegurarse de que todos los usuarios y administradores estén formados correctamente sobre cómo usar y mantener Contoso CipherGuard Sentinel X7."

# 9) "Documento e informes" / ": " ------------------------------------------
Replace-HeadingColon `
    "Documento e informes" `
    "Documentación e informes" `
    ": " `
    ":"

Replace-Simple `
    "Establezca hitos, asigne responsabilidades y realice un seguimiento del progreso." `
    "Establece hitos, asigna responsabilidades y realiza un seguimiento del progreso."

# 10) "Evaluación y mitigación" / " de riesgos:" ----------------------------
Replace-HeadingColon `
    "Evaluación y mitigación" `
    "Evaluación y mitigación de riesgos" `
    " de riesgos:" `
    ":"

# ---------------------------------------------------------------------------
# 11) "Escala de tiempo del proyecto." -> becomes a bold heading run plus a
#     separate, non-bold, colon-only run (mirrors the other headings).
# ---------------------------------------------------------------------------
$src = $d.Content
$srcOk = $src.Find.Execute("Pruebas y control de calidad", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $srcOk) {
    Write-Output "NOT FOUND: colon formatting donor paragraph"
} else {
    $srcColon = $d.Range($src.End, $src.End + 1)
    $savedFormatted = $srcColon.FormattedText

    $rng = $d.Content
    $ok = $rng.Find.Execute("Escala de tiempo del proyecto.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        Write-Output "NOT FOUND: Escala de tiempo del proyecto."
    } else {
        $start = $rng.Start
        $rng.Text = "Escala de tiempo del proyecto:"
        $newEnd = $start + ("Escala de tiempo del proyecto:").Length

        $wholeRng = $d.Range($start, $newEnd)
        $wholeRng.Font.Bold = 1
        $wholeRng.Font.BoldBi = 1

        $colonRng = $d.Range($newEnd - 1, $newEnd)
        $colonRng.FormattedText = $savedFormatted
    }
}

# ---------------------------------------------------------------------------
# Header: "con tecnología de IA," -> "Con tecnología de IA,"
# ---------------------------------------------------------------------------
$sec = $d.Sections(1)
$hdr = $sec.Headers.Item(1)
if ($hdr.Exists) {
    $hr = $hdr.Range
    $hOk = $hr.Find.Execute("con tecnología de IA,", $true, $false, $false, $false, $false, $true, 1, $false, "Con tecnología de IA,", 2)
    if (-not $hOk) {
        Write-Output "NOT FOUND: header text"
    }
}

Write-Output "ALL DONE"
